$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 17 ----
$ws.Range("A17").Value = 112175179
$ws.Range("B17").Value = 93388
$ws.Range("C17").Value = "Ovaliderad"
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 2180
$ws.Range("F17").Value = "Blåmossa"
$ws.Range("G17").Value = "Leucobryum glaucum"
$ws.Range("H17").Value = "(Hedw.) Ångstr."
$ws.Range("P17").Value = "Stora Körkroka (Stora Körkroka), Srm"
$ws.Range("Q17").Value = 693483.7626403375
$ws.Range("R17").Value = 6551529.771546691
$ws.Range("S17").Value = 25
$ws.Range("T17").Value = "Stockholm"
$ws.Range("U17").Value = "Haninge"
$ws.Range("V17").Value = "Södermanland"
$ws.Range("W17").Value = "Ornö"
$ws.Range("Y17").Value = "'2023-09-17"
$ws.Range("Z17").Value = "14:43"
$ws.Range("AA17").Value = "'2023-09-17"
$ws.Range("AB17").Value = "14:43"
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AW17").Value = "Måns Persson"
$ws.Range("AX17").Value = "Måns Persson"

# ---- Row 18 ----
$ws.Range("A18").Value = 112170170
$ws.Range("B18").Value = 56543
$ws.Range("C18").Value = "Ovaliderad"
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 103021
$ws.Range("F18").Value = "Talltita"
$ws.Range("G18").Value = "Poecile montanus"
$ws.Range("H18").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("P18").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q18").Value = 693556.9313844121
$ws.Range("R18").Value = 6551603.538506057
$ws.Range("S18").Value = 50
$ws.Range("T18").Value = "Stockholm"
$ws.Range("U18").Value = "Haninge"
$ws.Range("V18").Value = "Södermanland"
$ws.Range("W18").Value = "Ornö"
$ws.Range("Y18").Value = "'2023-09-17"
$ws.Range("Z18").Value = "13:51"
$ws.Range("AA18").Value = "'2023-09-17"
$ws.Range("AB18").Value = "13:51"
$ws.Range("AD18").Value = $false
$ws.Range("AE18").Value = $false
$ws.Range("AG18").Value = $false
$ws.Range("AW18").Value = "Klas Magnusson"
$ws.Range("AX18").Value = "Klas Magnusson, Måns Persson, Per Flodby"
